$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountOverview")

# New locator rows for the single-account balance verification feature.
# Values are entered in this specific order so the shared-string table
# is built up in the same sequence as the original authored workbook.
$ws.Range("A6").Value = "ELM_SingleAccountPane"
$ws.Range("C7").Value = "html/body/div[2]/div/div/div/div[3]/form/table/tbody/tr[1]/td[1]"
$ws.Range("C6").Value = "html/body/div[2]/div/div/div/div[3]/table[1]/tbody/tr[2]/td/table/tbody/tr/td[1]"
$ws.Range("A7").Value = "ELM_SingleAccountType"
$ws.Range("A8").Value = "ELM_SingleAccountBalance"
$ws.Range("C8").Value = "html/body/div[2]/div/div/div/div[3]/table[1]/tbody/tr[2]/td/table/tbody/tr/td[2]"
$ws.Range("B6").Value = "xpath"
$ws.Range("B7").Value = "xpath"
$ws.Range("B8").Value = "xpath"

# Widen columns A and C to fit the newly-added, longer element names/paths
# (matches the recorded best-fit widths of 23 and ~75.57 characters).
$ws.Columns.Item(1).ColumnWidth = 22.16666666666667
$ws.Columns.Item(3).ColumnWidth = 74.65

# Make AccountOverview the active sheet (moves tabSelected off MemberPayment)
# and land the selection on the last data row, like the recorded session.
$ws.Activate()
$null = $ws.Range("A8").Select()
